# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    populated with the per-fund holding detail for the new quarter (same
#    layout as the existing 2021-Qx sheets).
# 2. Insert a corresponding new top row into the "总计" sheet summarizing
#    the new quarter, pushing the existing rows down.
#
# NOTE: avoid "$var:" inside double-quoted strings (e.g. "D$r:G$r") -
# PowerShell parses the colon as a scope separator and mangles the text;
# use string concatenation for any dynamically-built A1 range address.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Seed the sheet with the "2021-Q4" layout (bold/bordered header row +
# bold/bordered index column) by copying it wholesale, then overwrite
# every cell with this quarter's real values below. A1 is skipped (it is
# blank on every sibling sheet - copying the whole A1:H4 block stamps an
# empty-but-present <c r="A1"/> node that shouldn't exist).
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1:H4").Copy($q1.Range("B1"))
$q4.Range("A2:A4").Copy($q1.Range("A2"))

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# The fund-code column (B) and the numeric-looking columns (D:G) are
# stored as text in this workbook; force text format before writing so
# leading zeros / literal decimal strings are preserved instead of being
# parsed into numbers.
$q1.Range("B2:B4").NumberFormat = "@"
$q1.Range("D2:G4").NumberFormat = "@"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "000593"
$q1.Range("C2").Value = "易方达标普全球高端消费品指数增强(QDII)-美元现汇"
$q1.Range("D2").Value = "1.93"
$q1.Range("E2").Value = "92.46"
$q1.Range("F2").Value = "5.84"
$q1.Range("G2").Value = "0.1127"
$q1.Range("H2").Value = 7

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "005676"
$q1.Range("C3").Value = "易方达标普全球高端消费品指数增强C(QDII) - 人民币"
$q1.Range("D3").Value = "1.93"
$q1.Range("E3").Value = "92.46"
$q1.Range("F3").Value = "5.84"
$q1.Range("G3").Value = "0.1127"
$q1.Range("H3").Value = 7

# Row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "118002"
$q1.Range("C4").Value = "易方达标普全球高端消费品指数增强A(QDII) - 人民币"
$q1.Range("D4").Value = "1.93"
$q1.Range("E4").Value = "92.46"
$q1.Range("F4").Value = "5.84"
$q1.Range("G4").Value = "0.1127"
$q1.Range("H4").Value = 7

# Drop the transient "quote prefix" style the NumberFormat/text coercion
# above stamped onto B2:B4 / D2:G4, so those cells end up unstyled like
# their siblings on the other quarter sheets.
$q1.Range("B2:B4").Style = "Normal"
$q1.Range("D2:G4").Style = "Normal"

$q1.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. New summary row in "总计", shifting the existing rows down by one
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Push existing data rows (2021-Q4, 2021-Q3, 2021-Q2) down one row,
# working bottom-to-top so nothing is clobbered before it is copied.
# Only B:D (date/count/value) are copied verbatim; column A is the row's
# 0-based index, which is re-numbered explicitly afterwards.
$tot.Range("B4:D4").Copy($tot.Range("B5"))
$tot.Range("B3:D3").Copy($tot.Range("B4"))
$tot.Range("B2:D2").Copy($tot.Range("B3"))

# A5 is brand new - copy A4's style (bold/bordered/centered) onto it,
# then fix up every index cell's value.
$tot.Range("A4").Copy($tot.Range("A5"))
$tot.Range("A5").Value = 3
$tot.Range("A4").Value = 2
$tot.Range("A3").Value = 1

# New first data row: 2022-Q1 summary
$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 3
$tot.Range("D2").Value = 0.34

$tot.Range("A1").Select()
